$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving the "text" nature of the
# original cell (the sheet stores numbers/percentages as literal text, e.g.
# "327.79" or "-0.98%"; without this, Excel's autodetection would convert
# them into real numbers/percentages and reformat the cell style).
function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# row -> column letter -> new value (only cells that actually change)
$changes = @{
    2  = @{ D = '327.79'; E = '-0.98%'; G = '19' }
    3  = @{ D = '43.82';  E = '5.25%';  G = '19' }
    4  = @{ D = '5.406';  E = '-4.87%'; G = '19' }
    5  = @{ E = '-3.07%'; G = '19' }
    6  = @{ D = '8.670';  E = '-1.50%'; G = '19' }
    7  = @{ D = '1.915';  E = '-4.58%'; G = '19' }
    8  = @{ D = '4.304';  E = '-3.60%'; G = '19' }
    9  = @{ D = '2.752';  E = '-5.44%'; G = '19' }
    10 = @{ D = '0.9435'; E = '2.02%';  G = '19' }
    11 = @{ D = '0.1188'; E = '-7.71%'; G = '19' }
    12 = @{ D = '0.1899'; E = '-4.10%'; G = '19' }
    13 = @{ D = '0.09584'; E = '1.50%'; G = '19' }
    14 = @{ D = '0.04173'; E = '8.62%'; G = '19' }
    15 = @{ D = '0.1070'; E = '0.84%';  G = '19' }
    16 = @{ D = '0.001279'; E = '-1.79%'; G = '19' }
    17 = @{ D = '0.005984'; E = '-1.96%'; G = '19' }
    18 = @{ B = 'LEO';                    C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo';                          D = '3.553';   E = '3.67%';  G = '19' }
    19 = @{ B = 'BitpandaEcosystemToken'; C = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best';          D = '0.3498';  E = '-0.15%'; G = '19' }
    20 = @{ B = 'MCDex';                  C = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb';                          D = '8.517';   E = '-2.94%'; G = '19' }
    21 = @{ B = 'ProBitToken';            C = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob';                    D = '0.1359';  E = '-0.85%'; G = '19' }
    22 = @{ B = 'ZBToken';                C = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb';                     D = '0.2604';  E = '3.83%';  G = '19' }
    23 = @{ B = 'CoinExToken';            C = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet';                D = '0.04367'; E = '-1.26%'; G = '19' }
    24 = @{ B = 'BitKan';                 C = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan';                     D = '0.001240'; E = '-2.65%'; G = '19' }
    25 = @{ B = 'HotbitToken';            C = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb';                 D = '0.004300'; E = '-1.31%'; G = '19' }
    26 = @{ G = '19' }
    27 = @{ E = '0.67%'; G = '19' }
    28 = @{ G = '19' }
    29 = @{ G = '19' }
    30 = @{ G = '19' }
    31 = @{ G = '19' }
    32 = @{ G = '19' }
    33 = @{ G = '19' }
    34 = @{ G = '19' }
    35 = @{ G = '19' }
    36 = @{ G = '19' }
    37 = @{ G = '19' }
    38 = @{ G = '19' }
    39 = @{ D = '0.02687'; E = '-5.01%'; G = '19' }
    40 = @{ D = '0.05500'; E = '-1.10%'; G = '19' }
    41 = @{ D = '0.007817'; E = '-1.63%'; G = '19' }
    42 = @{ E = '7.88%'; G = '19' }
    43 = @{ D = '0.1396'; E = '-2.68%'; G = '19' }
    44 = @{ D = '0.002125'; E = '3.06%'; G = '19' }
    45 = @{ D = '0.009624'; E = '-17.92%'; G = '19' }
    46 = @{ D = '0.00007106'; E = '2.66%'; G = '19' }
    47 = @{ E = '0.66%'; G = '19' }
    48 = @{ D = '0.003473'; E = '0.36%'; G = '19' }
    49 = @{ E = '0.34%'; G = '19' }
    50 = @{ E = '0.66%'; G = '19' }
    51 = @{ E = '0.66%'; G = '19' }
}

# Columns B and C hold non-numeric-looking text (coin names / URLs), so Excel's
# value autodetection leaves them as plain text already.
$plainTextCols = @('B', 'C')

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        if ($plainTextCols -contains $col) {
            $cell.Value = $cols[$col]
        } else {
            Set-TextValue $cell $cols[$col]
        }
    }
}
